$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9; existing rows 9-60 shift down to 10-61.
$ws.Rows("9:9").Insert()

# Fill the newly inserted row 9 with the new record.
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 45051
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 100112052
$ws.Range("G9").Value = "Albahaca"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 400
$ws.Range("L9").Value = 500
$ws.Range("M9").Value = 450
$ws.Range("N9").Value = "`$/paquete"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 450
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = "Hortaliza"

# Match the date cell style used by the rest of column D.
$ws.Range("D9").NumberFormat = $ws.Range("D10").NumberFormat
